$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$cv = $wb.Worksheets.Item("CONVERTION")
$tbl = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# Insert two new leave-card rows into Table1 ("UT(0-4-0)" before what is row
# 209, and "UT(0-4-2)" before what is the old VL(3-0-0) row), pushing every
# following row down by one each time.
# ---------------------------------------------------------------------------

# New row #1: before current row 209. Use the row that will be shifted into
# 210 as the formatting template (mirrors how Excel keeps the row's own look
# when a new row is inserted above it).
$ws.Rows(209).Insert()
$ws.Range("A210:K210").Copy()
$ws.Range("A209:K209").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B209").Value = "UT(0-4-0)"
$ws.Range("D209").Value = 0.5
$ws.Range("G209").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),""​"",Table1[[#This Row],[EARNED]])"
$ws.Range("G209").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# New row #2: before what is now row 212 (the old "VL(3-0-0)" row).
$ws.Rows(212).Insert()
$ws.Range("A213:K213").Copy()
$ws.Range("A212:K212").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("B212").Value = "UT(0-4-2)"
$ws.Range("D212").Value = 0.504
$ws.Range("G212").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# Grow the table definition to cover the two new rows (A8:K345 -> A8:K347).
$tbl.Resize($ws.Range("A8:K347"))

# ---------------------------------------------------------------------------
# CONVERTION sheet: update the undertime lookup used for the new entries.
# ---------------------------------------------------------------------------
$cv.Range("E3").Value = 4
$cv.Range("F3").ClearContents()

# ---------------------------------------------------------------------------
# Leave the selection where the author left it after making the edit.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("F211").Select()

$wb.Application.Calculate()
